$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D to make room for the new quarter
# columns (Dec-2018 and Sep-2018), shifting all existing quarterly data right.
$ws.Range("D:E").Insert()

# The newly inserted columns inherit the formatting of column C (text).
# Copy the number/date formatting from column F (the old column D, now shifted)
# onto the new D:E columns so the new quarters match the existing style.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the two new quarter columns (D = Dec-2018, E = Sep-2018) with data
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 651000
$ws.Range("E8").Value = 633200
$ws.Range("D9").Value = 514000
$ws.Range("E9").Value = 482000
$ws.Range("D10").Value = 137000
$ws.Range("E10").Value = 151200
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 100
$ws.Range("E14").Value = 600
$ws.Range("D15").Value = 27500
$ws.Range("E15").Value = 31800
$ws.Range("D17").Value = 564500
$ws.Range("E17").Value = 550900
$ws.Range("D18").Value = 86500
$ws.Range("E18").Value = 82300
$ws.Range("D20").Value = 1100
$ws.Range("E20").Value = 800
$ws.Range("D21").Value = 115100
$ws.Range("E21").Value = 114900
$ws.Range("D22").Value = 4800
$ws.Range("E22").Value = 5200
$ws.Range("D23").Value = 82700
$ws.Range("E23").Value = 78000
$ws.Range("D24").Value = -3400
$ws.Range("E24").Value = -45200
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 86100
$ws.Range("E26").Value = 123200
$ws.Range("D27").Value = 86100
$ws.Range("E27").Value = 123200
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = -1100
$ws.Range("E32").Value = -800
$ws.Range("D33").Value = 86100
$ws.Range("E33").Value = 123200
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 86100
$ws.Range("E35").Value = 123200
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 264900
$ws.Range("E41").Value = 245700
$ws.Range("D42").Value = 162800
$ws.Range("E42").Value = 162500
$ws.Range("D43").Value = 249800
$ws.Range("E43").Value = 210300
$ws.Range("D44").Value = 125500
$ws.Range("E44").Value = 163900
$ws.Range("D45").Value = 75700
$ws.Range("E45").Value = 95000
$ws.Range("D46").Value = 878800
$ws.Range("E46").Value = 877400
$ws.Range("D47").Value = 104700
$ws.Range("E47").Value = 105300
$ws.Range("D48").Value = 834800
$ws.Range("E48").Value = 919600
$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 68800
$ws.Range("E52").Value = 95300
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 1887100
$ws.Range("E54").Value = 1997600
$ws.Range("D57").Value = 128000
$ws.Range("E57").Value = 119600
$ws.Range("D58").Value = 17800
$ws.Range("E58").Value = 11500
$ws.Range("D59").Value = 183500
$ws.Range("E59").Value = 203700
$ws.Range("D60").Value = 329300
$ws.Range("E60").Value = 334800
$ws.Range("D61").Value = 300200
$ws.Range("E61").Value = 302800
$ws.Range("D62").Value = 552700
$ws.Range("E62").Value = 672900
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 1182200
$ws.Range("E66").Value = 1310500
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 527700
$ws.Range("E72").Value = 449100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 704800
$ws.Range("E76").Value = 687100
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 86100
$ws.Range("E81").Value = 123200
$ws.Range("D83").Value = 27500
$ws.Range("E83").Value = 31800
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 151800
$ws.Range("E89").Value = 121100
$ws.Range("D91").Value = -39500
$ws.Range("E91").Value = -25700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -39700
$ws.Range("E94").Value = -28900
$ws.Range("D96").Value = -7300
$ws.Range("E96").Value = -7600
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -92900
$ws.Range("E100").Value = -88100
$ws.Range("D101").Value = 0
$ws.Range("E101").Value = 0
$ws.Range("D102").Value = 19300
$ws.Range("E102").Value = 4100

# Restated / corrected historical figures for a few earlier quarters
$ws.Range("H9").Value = 450700
$ws.Range("I9").Value = 494400
$ws.Range("H10").Value = 109500
$ws.Range("I10").Value = 119100
$ws.Range("H17").Value = 508600
$ws.Range("I17").Value = 540800
$ws.Range("H18").Value = 51600
$ws.Range("I18").Value = 72700
$ws.Range("H20").Value = 400
$ws.Range("I20").Value = -100
$ws.Range("H24").Value = -367200
$ws.Range("H26").Value = 413700
$ws.Range("H27").Value = 413700
$ws.Range("H29").Value = -332400
$ws.Range("I29").Value = 0
$ws.Range("H32").Value = -400
$ws.Range("I32").Value = 100
$ws.Range("F91").Value = -20600
$ws.Range("H91").Value = -28700
$ws.Range("I91").Value = -13600
$ws.Range("J91").Value = -11000
